$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1ST Q 2024")

# --- Row 12: TOLENTINO, ABRAHAM N / CITY MAYOR -> ITALY, MARCH 08 - 12, 2024 ---
$ws.Range("A12").Value = "TOLENTINO, ABRAHAM N"
$ws.Range("B12").Value = "CITY MAYOR"
$ws.Range("C12").Value = "ITALY"
$ws.Range("D12").Value = "MARCH 08 - 12, 2024"
$ws.Range("E12").Value = "PERSONAL"
$ws.Range("H12").Value = "/"
$ws.Range("I12").Value = 1

# --- Row 13: TOLENTINO, ABRAHAM N / CITY MAYOR -> THAILAND, MARCH 27 - 31, 2024 ---
$ws.Range("A13").Value = "TOLENTINO, ABRAHAM N"
$ws.Range("B13").Value = "CITY MAYOR"
$ws.Range("C13").Value = "THAILAND"
$ws.Range("D13").Value = "MARCH 27 - 31, 2024"
$ws.Range("E13").Value = "PERSONAL"
$ws.Range("H13").Value = "/"
$ws.Range("I13").Value = 1

# Match the "/" cell formatting (general number format w/ quote-prefix) used
# by the other filled-in rows (e.g. H11) instead of the blank template's
# text-number-format. Do this AFTER the values are set, since pasting over a
# value resets the formatting.
$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The newly-filled rows wrap text like the other data rows above them, which
# grows the row height the same way.
$ws.Rows.Item(12).RowHeight = 28.8
$ws.Rows.Item(13).RowHeight = 28.8

# Restore the cursor/selection to where it ended up after this edit.
[void]$ws.Range("H10:I13").Select()
